# Refresh the cryptos price-list snapshot (Price / Volume(1h) columns) to the
# latest scraped values. Numeric-looking prices are entered with a leading
# apostrophe so Excel keeps them as literal text (matching how this sheet
# already stores "Price" as text, e.g. "26.656.58", "19.70", "0.0172") instead
# of silently re-typing them as numbers and dropping trailing zeros / digit
# grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $ws.Range($Address).Value = "'" + $Text
}

# Row 2 - Bitcoin
Set-TextValue 'D2' '26.656.58'

# Row 3 - Ethereum
Set-TextValue 'D3' '1.592.42'
$ws.Range('E3').Value = '  -2.54%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.03%  '

# Row 5 - BNB
Set-TextValue 'D5' '211.14'
$ws.Range('E5').Value = '  -2.42%  '

# Row 6 - XRP
$ws.Range('E6').Value = '  -2.12%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  -0.02%  '

# Row 8 - Cardano
$ws.Range('E8').Value = '  -3.22%  '

# Row 9 - Dogecoin
$ws.Range('E9').Value = '  -1.70%  '

# Row 10 - Solana
Set-TextValue 'D10' '19.70'
$ws.Range('E10').Value = '  -3.12%  '

# Row 11 - TRON
$ws.Range('E11').Value = '  -1.83%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range('E12').Value = '  -2.64%  '

# Row 13 - WrappedEther
Set-TextValue 'D13' '1.592.66'
$ws.Range('E13').Value = '  -2.54%  '

# Row 15 - Polygon
$ws.Range('E15').Value = '  -3.41%  '

# Row 16 - Litecoin
Set-TextValue 'D16' '64.79'
$ws.Range('E16').Value = '  -0.70%  '

# Row 17 - WrappedBTC
Set-TextValue 'D17' '26.651.83'
$ws.Range('E17').Value = '  -1.91%  '

# Row 18 - ShibaInu
Set-TextValue 'D18' '0.0₃0727'

# Row 19 - BitcoinCash
Set-TextValue 'D19' '208.30'
$ws.Range('E19').Value = '  -4.46%  '

# Row 20 - Dai
$ws.Range('E20').Value = '  -0.13%  '

# Row 21 - Chainlink
Set-TextValue 'D21' '6.77'
$ws.Range('E21').Value = '  -2.56%  '

# Row 22 - Uniswap
Set-TextValue 'D22' '4.26'
$ws.Range('E22').Value = '  -3.36%  '

# Row 23 - Toncoin
$ws.Range('E23').Value = '  -1.68%  '

# Row 24 - Avalanche
Set-TextValue 'D24' '8.91'
$ws.Range('E24').Value = '  -1.72%  '

# Row 25 - Monero
Set-TextValue 'D25' '147.10'
$ws.Range('E25').Value = '  -0.66%  '

# Row 26 - BinanceUSD
$ws.Range('E26').Value = '  +0.00%  '

# Row 27 - Cosmos
Set-TextValue 'D27' '7.28'

# Row 28 - Stellar
$ws.Range('E28').Value = '  -3.91%  '

# Row 29 - EthereumClassic
Set-TextValue 'D29' '15.32'
$ws.Range('E29').Value = '  -2.31%  '

# Row 30 - Hedera
$ws.Range('E30').Value = '  -0.54%  '

# Row 31 - PancakeSwap
$ws.Range('E31').Value = '  -2.07%  '

# Row 32 - Filecoin
$ws.Range('E32').Value = '  -4.50%  '

# Row 33 - ImmutableX
Set-TextValue 'D33' '0.653'
$ws.Range('E33').Value = '  +18.83%  '

# Row 34 - InternetComputer(DFINITY)
Set-TextValue 'D34' '2.90'
$ws.Range('E34').Value = '  -3.52%  '

# Row 35 - Maker
Set-TextValue 'D35' '1.316.00'
$ws.Range('E35').Value = '  -2.05%  '

# Row 36 - LidoDAOToken
$ws.Range('E36').Value = '  -4.87%  '

# Row 37 - HuobiToken
$ws.Range('E37').Value = '  -2.13%  '

# Row 38 - VeChain
Set-TextValue 'D38' '0.0172'
$ws.Range('E38').Value = '  -2.51%  '

# Row 39 - ARBITRUM
$ws.Range('E39').Value = '  -2.85%  '

# Row 42 - TrustWalletToken
$ws.Range('E42').Value = '  -1.74%  '

# Row 43 - MXToken
$ws.Range('E43').Value = '  -3.59%  '

# Row 44 - Aave
Set-TextValue 'D44' '63.53'
$ws.Range('E44').Value = '  -1.86%  '

# Row 45 - RocketPoolETH
Set-TextValue 'D45' '1.727.61'
$ws.Range('E45').Value = '  -2.62%  '

# Row 46 - Quant
Set-TextValue 'D46' '89.94'
$ws.Range('E46').Value = '  -1.00%  '

# Row 47 - RenderToken
$ws.Range('E47').Value = '  -0.91%  '

# Row 48 - WEMIXToken
$ws.Range('E48').Value = '  +3.20%  '

# Row 49 - Cronos
Set-TextValue 'D49' '0.0509'
$ws.Range('E49').Value = '  -0.93%  '

# Row 50 - Algorand
Set-TextValue 'D50' '0.0978'
$ws.Range('E50').Value = '  -1.30%  '

# Row 51 - EnergySwap
Set-TextValue 'D51' '7.50'
$ws.Range('E51').Value = '  -1.36%  '
